$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 0.93
$ws.Range("C5").Value = 0.9399999999999999
$ws.Range("E5").Value = 0.61
$ws.Range("F5").Value = 0.46
$ws.Range("H5").Value = 0.63
$ws.Range("I5").Value = 0.46
$ws.Range("K5").Value = 0.75
$ws.Range("L5").Value = 0.61
$ws.Range("N5").Value = 0.85
$ws.Range("O5").Value = 0.77
